$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.383.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.306.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.305.69"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.85"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E11").Value = "  -4.50%  "
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.870.59"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.302.92"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000166"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.393.76"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.537"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.435.75"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000103"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -8.60%  "
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.69%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.57"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.64"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.25"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.87"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.70"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.336.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -14.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0732"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.751"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.14"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.351.30"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.45%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.40"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.31%  "
